# Update the admission form with the new patient's data
# (per commit message: "para agregar codigo de barras" - data refreshed
# for a new patient record on this report).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Apellidos / Nombres / No. Expediente (row 6)
$ws.Range("A6").Value = "VETORAZZI"
$ws.Range("C6").Value = "BARRIOS"
$ws.Range("E6").Value = "ELVA"
$ws.Range("G6").Value = "VIOLETA"
$ws.Range("I6").Value = "/201773425"

# Direccion actual: Calle o lugar / Municipio / Departamento (row 8)
$ws.Range("A8").Value = "6 AV 1-60"
$ws.Range("D8").Value = "Z. 19 LA FLORIDA"
$ws.Range("F8").Value = "GUATEMALA"
$ws.Range("H8").Value = "GUATEMALA"

# Fecha de nacimiento / Edad / Lugar de nacimiento / Sexo (row 12)
$ws.Range("A12").Value = "1955-07-26"
$ws.Range("F12").Value = "62"
$ws.Range("H12").Value = "SUCHITEPEQUEZ"
$ws.Range("J12").Value = "Femenino"

# Estado civil / Ocupacion / Nacionalidad / No. de Cedula (row 14)
$ws.Range("A14").Value = "Casado"
$ws.Range("D14").Value = "AMA DE CASA"
$ws.Range("F14").Value = "GUATEMALTECA"
$ws.Range("H14").Value = "NO PRESENTA"

# Direccion si difiere a la indicada (row 16)
$ws.Range("A16").Value = "MARCO TULIO ALVAREZ"

# Nombre del Padre / Nombre de la Madre (row 18)
$ws.Range("A18").Value = "JUAN PABLO VETORAZZI"
$ws.Range("F18").Value = "MARIA EUGENIA BARRIOS"

# En caso de emergencia notificar a / Parentesco / Direccion / Telefono (row 20)
$ws.Range("A20").Value = "JAQUELIN ALVAREZ"
$ws.Range("F20").Value = "HIJA"
$ws.Range("H20").Value = ""
$ws.Range("J20").Value = "3406-8429"

# Fecha de Ingreso / Hora (row 24)
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "13:4:40"
